$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 200, shifting existing rows 200-272 down to 201-273.
$ws.Rows.Item(200).Insert()

# Populate the newly inserted row 200 with a new weekly record. The
# descriptive columns (A, B, C, E, F, G, H, I, R) carry the same constant
# values used throughout this subset (market/region/category/quality/unit).
$ws.Cells.Item(200, 1).Value = 9
$ws.Cells.Item(200, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(200, 3).Value = "Metropolitana"
$ws.Cells.Item(200, 4).Value = 44704
$ws.Cells.Item(200, 5).Value = 13
$ws.Cells.Item(200, 6).Value = 100112001
$ws.Cells.Item(200, 7).Value = "Berenjena"
$ws.Cells.Item(200, 8).Value = "Sin especificar"
$ws.Cells.Item(200, 9).Value = "Primera"
$ws.Cells.Item(200, 10).Value = 190
$ws.Cells.Item(200, 11).Value = 6000
$ws.Cells.Item(200, 12).Value = 7000
$ws.Cells.Item(200, 13).Value = 6526
$ws.Cells.Item(200, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(200, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(200, 16).Value = 131
$ws.Cells.Item(200, 17).Value = 50
$ws.Cells.Item(200, 18).Value = "Hortaliza"

# Match the date number format used by the rest of column D.
$ws.Cells.Item(200, 4).NumberFormat = $ws.Cells.Item(201, 4).NumberFormat
